$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.44 = 30215.3 pesos`n✅ 30215.3 pesos = 7.4 = 966.55 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 134.369
$wsTasas.Range("O10").Value = 4060
$wsTasas.Range("N12").Value = 4083
$wsTasas.Range("O12").Value = 130.61
